$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "FirstName"
$ws.Range("F2").Value = "Ravi"
$ws.Range("G1").Value = "LastName"
$ws.Range("G2").Value = "Kumar"

$ws.Range("E3").Select()
